# The source data is a list of unique transit-station names (one per row in
# column A). This edit corrects three station names by querying the ODsay
# transit API (per the commit message "feat: call odsay api"):
#
#   - row 28  "화전"        -> "한국항공대"   (renamed in place)
#   - row 194 "도라산"      -> removed entirely (row deleted, no replacement)
#   - row 410 "뚝섬유원지"  -> "자양"          (renamed in place)
#
# Net effect: 446 rows -> 445 rows (one row removed), dimension A1:A446 -> A1:A445.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the station at row 28 ("화전" -> "한국항공대").
$ws.Range("A28").Value = "한국항공대"

# 2) Delete the row that held "도라산" outright (row 194 in the original
#    layout); everything below shifts up by one row.
$ws.Rows.Item(194).Delete()

# 3) Rename the station that is now at row 409 (originally row 410,
#    "뚝섬유원지", shifted up by the deletion above) to "자양".
$ws.Range("A409").Value = "자양"

# Reflect the author's final selection (cell A410, matching the saved
# sheetView) so the workbook opens with the same active cell.
$ws.Range("A410").Select()
